{"js": "// Diary entry for \"2023\u5e743\u67089\u65e5\u661f\u671f\u56db \u6674\" needs to stay as its own line,\n// and a NEW line for \"2023\u5e743\u670810\u65e5\u661f\u671f\u4e94 \u6674\" must follow it (the old\n// paragraph keeps its own formatting, but its text changes to the 3/10\n// entry) per the commit \"add words in March 10th\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"2023\u5e743\u67089\u65e5\u661f\u671f\u56db \u6674\";\n\nlet dateParagraph = null;\nlet priorParagraph = null;\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(marker) !== -1) {\n    dateParagraph = items[i];\n    priorParagraph = i > 0 ? items[i - 1] : null;\n    break;\n  }\n}\n\nif (!dateParagraph) {\n  throw new Error(\"Could not find the '\" + marker + \"' paragraph.\");\n}\n\n// Insert a fresh paragraph carrying the original March 9th text right\n// after the paragraph that preceded the date line (this reproduces the\n// paragraph-mark formatting Word itself would inherit from pressing Enter\n// at the end of that preceding line).\nif (priorParagraph) {\n  priorParagraph.insertParagraph(marker, \"After\");\n} else {\n  dateParagraph.insertParagraph(marker, \"Before\");\n}\n\n// The paragraph that used to hold the 3/9 entry now becomes the 3/10 entry.\ndateParagraph.insertText(\"2023\u5e743\u670810\u65e5\u661f\u671f\u4e94 \u6674\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Diary entry for \"2023\u5e743\u67089\u65e5\u661f\u671f\u56db \u6674\" needs to stay as its own line,\n# and a NEW line for \"2023\u5e743\u670810\u65e5\u661f\u671f\u4e94 \u6674\" must follow it (the old\n# paragraph keeps its own position but gets the new text) per the commit\n# \"add words in March 10th\".\n$d = $word.ActiveDocument\n$marker = \"2023\u5e743\u67089\u65e5\u661f\u671f\u56db \u6674\"\n\n$searchRange = $d.Content\n$found = $searchRange.Find.Execute($marker)\nif (-not $found) {\n    throw \"Could not find the '$marker' paragraph.\"\n}\n\n# Expand the hit to the full paragraph so we know exactly where it starts.\n$searchRange.Expand(4) | Out-Null   # 4 = wdParagraph\n$origStart = $searchRange.Start\n\n# Figure out which 1-based Paragraphs() index this is, so we can address\n# both the new and the original paragraph afterwards.\n$paraIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Start -eq $origStart) {\n        $paraIndex = $i\n        break\n    }\n}\nif ($paraIndex -eq 0) {\n    throw \"Could not resolve the paragraph index for '$marker'.\"\n}\n\nif ($origStart -gt 0) {\n    # Insert the new paragraph mark right after the END of the PRECEDING\n    # paragraph (rather than right before the target paragraph) so the new\n    # paragraph mark inherits that preceding paragraph's run formatting\n    # (matches what Word does when you place the caret at the end of the\n    # previous line and press Enter).\n    $anchor = $d.Range($origStart - 1, $origStart - 1)\n    $anchor.InsertParagraphAfter() | Out-Null\n} else {\n    # Target paragraph is the very first paragraph in the document - fall\n    # back to inserting directly before it.\n    $searchRange.InsertParagraphBefore() | Out-Null\n}\n\n# The freshly inserted paragraph takes the original March 9th text \u2026\n$d.Paragraphs.Item($paraIndex).Range.Text = \"2023\u5e743\u67089\u65e5\u661f\u671f\u56db \u6674\"\n# \u2026 and the paragraph that used to hold it becomes the March 10th entry.\n$d.Paragraphs.Item($paraIndex + 1).Range.Text = \"2023\u5e743\u670810\u65e5\u661f\u671f\u4e94 \u6674\"\n"}
